# Auto-generated Excel COM-interop script applying the scheduled-runner
# market-data refresh described in the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 250001330
$ws.Range("I18").Value = 280
$ws.Range("J18").Value = 333335000
$ws.Range("K18").Value = 280
$ws.Range("L18").Value = 333335000
$ws.Range("M18").Value = 4
$ws.Range("N18").Value = -333335568
$ws.Range("H19").Value = 1062.7142
$ws.Range("J19").Value = 1128
$ws.Range("L19").Value = 1128
$ws.Range("N19").Value = -1478
$ws.Range("H33").Value = 211.125
$ws.Range("I33").Value = 211.125
$ws.Range("K33").Value = 211.125
$ws.Range("M33").Value = 17.875
$ws.Range("H40").Value = 6875.25
$ws.Range("J40").Value = 7333.6665
$ws.Range("L40").Value = 7333.6665
$ws.Range("N40").Value = -7683.6665
$ws.Range("H58").Value = 681.125
$ws.Range("J58").Value = 5000
$ws.Range("L58").Value = 15000
$ws.Range("N58").Value = -15300
$ws.Range("H98").Value = 3896.2942
$ws.Range("J98").Value = 1553
$ws.Range("L98").Value = 1553
$ws.Range("N98").Value = -4549
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("H122").Value = 3896.2942
$ws.Range("J122").Value = 1553
$ws.Range("L122").Value = 4659
$ws.Range("N122").Value = -9559
$ws.Range("H132").Value = 5761.28
$ws.Range("I132").Value = 6142.7393
$ws.Range("K132").Value = 18428.2179
$ws.Range("M132").Value = -15898.2179
$ws.Range("H137").Value = 1989.0588
$ws.Range("I137").Value = 2048
$ws.Range("J137").Value = 1797.5
$ws.Range("K137").Value = 6144
$ws.Range("L137").Value = 5392.5
$ws.Range("M137").Value = -3594
$ws.Range("N137").Value = -10492.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 153.28572
$ws.Range("I4").Value = 129.83333
$ws.Range("J4").Value = 294
$ws.Range("K4").Value = 129.83333
$ws.Range("L4").Value = 294
$ws.Range("M4").Value = -13.83332999999999
$ws.Range("N4").Value = -526
$ws.Range("H61").Value = 2898.8096
$ws.Range("I61").Value = 2348.2856
$ws.Range("K61").Value = 2348.2856
$ws.Range("M61").Value = -2136.2856
$ws.Range("H70").Value = 89999
$ws.Range("J70").Value = 89999
$ws.Range("L70").Value = 89999
$ws.Range("N70").Value = -90539
$ws.Range("H73").Value = 89999
$ws.Range("J73").Value = 89999
$ws.Range("L73").Value = 89999
$ws.Range("N73").Value = -91871
$ws.Range("H74").Value = 1163.3096
$ws.Range("I74").Value = 704.2646999999999
$ws.Range("K74").Value = 704.2646999999999
$ws.Range("M74").Value = 169.7353000000001
$ws.Range("H77").Value = 1163.3096
$ws.Range("I77").Value = 704.2646999999999
$ws.Range("K77").Value = 3521.3235
$ws.Range("M77").Value = 846.6765000000005
$ws.Range("H95").Value = 8666.333000000001
$ws.Range("I95").Value = 7000
$ws.Range("J95").Value = 9499.5
$ws.Range("K95").Value = 7000
$ws.Range("L95").Value = 9499.5
$ws.Range("M95").Value = -4254
$ws.Range("N95").Value = -14991.5
$ws.Range("H136").Value = 2898.8096
$ws.Range("I136").Value = 2348.2856
$ws.Range("K136").Value = 7044.8568
$ws.Range("M136").Value = -4494.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 21669308
$ws.Range("I105").Value = 1252787.4
$ws.Range("J105").Value = 62502348
$ws.Range("K105").Value = 1252787.4
$ws.Range("L105").Value = 62502348
$ws.Range("M105").Value = -1251040.4
$ws.Range("N105").Value = -62505842
$ws.Range("H134").Value = 1889.4615
$ws.Range("I134").Value = 1201.5
$ws.Range("K134").Value = 3604.5
$ws.Range("M134").Value = -1069.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 649.75
$ws.Range("I22").Value = 533
$ws.Range("K22").Value = 533
$ws.Range("M22").Value = -183
$ws.Range("H31").Value = 6582211
$ws.Range("I31").Value = 2322.923
$ws.Range("K31").Value = 2322.923
$ws.Range("M31").Value = -2027.923
$ws.Range("H34").Value = 6582211
$ws.Range("I34").Value = 2322.923
$ws.Range("K34").Value = 2322.923
$ws.Range("M34").Value = -2120.923
$ws.Range("H58").Value = 1894.6471
$ws.Range("I58").Value = 1246.2222
$ws.Range("K58").Value = 1246.2222
$ws.Range("M58").Value = -1043.2222
$ws.Range("H74").Value = 59437.332
$ws.Range("J74").Value = 59437.332
$ws.Range("L74").Value = 59437.332
$ws.Range("N74").Value = -61185.332
$ws.Range("H77").Value = 59437.332
$ws.Range("J77").Value = 59437.332
$ws.Range("L77").Value = 178311.996
$ws.Range("N77").Value = -187047.996
$ws.Range("H94").Value = 479.72726
$ws.Range("I94").Value = 514
$ws.Range("K94").Value = 514
$ws.Range("M94").Value = -63
$ws.Range("H122").Value = 2222.52
$ws.Range("J122").Value = 4058.8
$ws.Range("L122").Value = 12176.4
$ws.Range("N122").Value = -17076.4
$ws.Range("H132").Value = 3904.9143
$ws.Range("I132").Value = 3319.7083
$ws.Range("K132").Value = 9959.124899999999
$ws.Range("M132").Value = -7429.124899999999
$ws.Range("H134").Value = 3462.9023
$ws.Range("I134").Value = 3477.1562
$ws.Range("K134").Value = 10431.4686
$ws.Range("M134").Value = -7896.4686
$ws.Range("H136").Value = 1894.6471
$ws.Range("I136").Value = 1246.2222
$ws.Range("K136").Value = 3738.6666
$ws.Range("M136").Value = -1188.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 677.8
$ws.Range("I5").Value = 677.8
$ws.Range("K5").Value = 2033.4
$ws.Range("M5").Value = -1921.4
$ws.Range("H120").Value = 500
$ws.Range("I120").Value = 500
$ws.Range("K120").Value = 1500
$ws.Range("M120").Value = 3338
$ws.Range("H131").Value = 2684287.8
$ws.Range("J131").Value = 3923486.8
$ws.Range("L131").Value = 11770460.4
$ws.Range("N131").Value = -11780540.4
$ws.Range("H133").Value = 2330
$ws.Range("I133").Value = 995
$ws.Range("K133").Value = 2985
$ws.Range("M133").Value = 2075
$ws.Range("H134").Value = 4724.5
$ws.Range("I134").Value = 4265.6665
$ws.Range("J134").Value = 4999.8
$ws.Range("K134").Value = 12796.9995
$ws.Range("L134").Value = 14999.4
$ws.Range("M134").Value = -7726.999500000002
$ws.Range("N134").Value = -25139.4
$ws.Range("H135").Value = 677.8
$ws.Range("I135").Value = 677.8
$ws.Range("K135").Value = 6100.2
$ws.Range("M135").Value = -3565.2
$ws.Range("H138").Value = 15000
$ws.Range("I138").Value = 10000
$ws.Range("K138").Value = 30000
$ws.Range("M138").Value = -24860
$ws.Range("H139").Value = 5260.4595
$ws.Range("I139").Value = 3633.5625
$ws.Range("J139").Value = 6500
$ws.Range("K139").Value = 10900.6875
$ws.Range("L139").Value = 19500
$ws.Range("M139").Value = -5760.6875
$ws.Range("N139").Value = -29780

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 7000
$ws.Range("I43").Value = 7000
$ws.Range("K43").Value = 7000
$ws.Range("M43").Value = -6849
$ws.Range("H107").Value = 1375.421
$ws.Range("I107").Value = 440.77777
$ws.Range("K107").Value = 440.77777
$ws.Range("M107").Value = 1479.22223
$ws.Range("H132").Value = 2297.6765
$ws.Range("I132").Value = 2230.2173
$ws.Range("J132").Value = 2438.7273
$ws.Range("K132").Value = 6690.651899999999
$ws.Range("L132").Value = 7316.1819
$ws.Range("M132").Value = -4160.651899999999
$ws.Range("N132").Value = -12376.1819

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1080.5714
$ws.Range("I16").Value = 1073
$ws.Range("J16").Value = 1099.5
$ws.Range("K16").Value = 1073
$ws.Range("L16").Value = 1099.5
$ws.Range("M16").Value = -903
$ws.Range("N16").Value = -1439.5
$ws.Range("H64").Value = 23016
$ws.Range("J64").Value = 23016
$ws.Range("L64").Value = 23016
$ws.Range("N64").Value = -23466
$ws.Range("H67").Value = 23016
$ws.Range("J67").Value = 23016
$ws.Range("L67").Value = 23016
$ws.Range("N67").Value = -24576
$ws.Range("H132").Value = 4430.778
$ws.Range("I132").Value = 4742.364
$ws.Range("J132").Value = 3941.1428
$ws.Range("K132").Value = 14227.092
$ws.Range("L132").Value = 11823.4284
$ws.Range("M132").Value = -11697.092
$ws.Range("N132").Value = -16883.4284
$ws.Range("H136").Value = 4789.9287
$ws.Range("I136").Value = 4327.909
$ws.Range("J136").Value = 6484
$ws.Range("K136").Value = 12983.727
$ws.Range("L136").Value = 19452
$ws.Range("M136").Value = -10433.727

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 10871456
$ws.Range("I122").Value = 1798.9474
$ws.Range("K122").Value = 5396.8422
$ws.Range("M122").Value = -2946.8422
$ws.Range("H132").Value = 2537.3845
$ws.Range("I132").Value = 2602.92
$ws.Range("J132").Value = 899
$ws.Range("K132").Value = 7808.76
$ws.Range("L132").Value = 2697
$ws.Range("M132").Value = -5278.76
$ws.Range("N132").Value = -7757

